$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Insert row at 13 (before old row12=Misc header at that point... let's just test row 9 for example, where row8 above is blank s6, row9 below also blank s6)
$ws.Rows.Item(9).Insert()
